# Generate Report for Handback
# - Mark the 94ba0302... file's status as a failed handback transform
#   (Overview!E3/F3, zh-cn!C3, de-de!C3).
# - Record the handback/handoff filename-mismatch error detail for that
#   file in the per-locale sheets' "Error Detail" column (P3).
# - Widen the "Error Detail" column so the new message is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status of the 94ba0302-... row changes from "Ready for handoff" to
# "Handback transform failed" everywhere it is shown.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# New "Error Detail" messages for the handback filename mismatch.
$wsZhCn.Range("P3").Value = "Handback file name: ijwtizsi.fze is different with handoff file name: 94ba0302-8744-4d0c-ba8a-6ce35f1e8b60.fe68e2a71c998b673fef630e578798e9128d8693.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: ijwtizsi.fze is different with handoff file name: 94ba0302-8744-4d0c-ba8a-6ce35f1e8b60.fe68e2a71c998b673fef630e578798e9128d8693.de-de."

# Widen the "Error Detail" column (P / column 16) to fit the new text.
# ColumnWidth of 39.1 round-trips through Excel's pixel-quantised storage
# to a saved <col width="40">, matching the target width exactly.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1
